$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Output "WARNING: could not find '$old'"
    }
}

Replace-Text "2025-04-22 Tuesday" "2025-04-23 Wednesday"

Replace-Text "69×29=" "48×12="
Replace-Text "55×29=" "92×35="
Replace-Text "15×91=" "79×87="
Replace-Text "35×20=" "74×71="
Replace-Text "38×41=" "95×89="
Replace-Text "75×32=" "79×41="
Replace-Text "31×30=" "60×88="
Replace-Text "59×94=" "55×95="
Replace-Text "77×23=" "25×71="
Replace-Text "83×42=" "43×44="
Replace-Text "12×69=" "28×38="
Replace-Text "65×44=" "74×36="
Replace-Text "61×56=" "99×84="
Replace-Text "38×42=" "87×24="
Replace-Text "62×49=" "58×82="
Replace-Text "20×74=" "32×33="
Replace-Text "96×55=" "37×20="
Replace-Text "97×15=" "89×14="
Replace-Text "50×47=" "20×78="
Replace-Text "91×29=" "60×52="
Replace-Text "12×96=" "16×14="
Replace-Text "31×74=" "15×38="
Replace-Text "62×56=" "52×66="
Replace-Text "71×53=" "87×77="
Replace-Text "53×39=" "51×55="
